# ES-272 - updated experiments-report-template.xlsx
#
# The report template's "Email заявки" / ${experiment.email} column is
# replaced by a "Пользователь" / ${experiment.createdBy} column, and the
# sheet's saved view (scroll position / selection) is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Replace the "Email заявки" header / placeholder column (E6/E7) -------
$ws.Range("E6").Value = "Пользователь"
$ws.Range("E7").Value = "`${experiment.createdBy}"

# --- Update the saved sheet view: scrolled to column B, selection on B7 ---
$ws.Activate()
$ws.Range("B7").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
